# Edit script: fix STD/ABSA code finding procedures
# 1. Swap the D1/E1 header labels (CREDIT/DEBIT -> DEBIT/CREDIT)
# 2. Clean up DESCRIPTION_CODE (column B) values:
#    - rows without a CODE1 (column C) value: trim only trailing whitespace
#    - rows with a CODE1 value: trim, collapse internal whitespace, replace
#      hyphens with spaces, and upper-case the text

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Swap header labels in D1 / E1 ---
$ws.Cells.Item(1, 4).Value = "DEBIT"
$ws.Cells.Item(1, 5).Value = "CREDIT"

# --- 2. Normalize DESCRIPTION_CODE column (B) ---
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($row = 2; $row -le $lastRow; $row++) {
    $bCell = $ws.Cells.Item($row, 2)
    $cCell = $ws.Cells.Item($row, 3)

    $bValue = $bCell.Value()
    if ($null -eq $bValue) { continue }
    $bValue = [string]$bValue

    $cValue = $cCell.Value()

    if ($null -eq $cValue -or [string]$cValue -eq "") {
        $newValue = $bValue.TrimEnd()
    } else {
        $newValue = $bValue.Trim()
        $newValue = $newValue.Replace("-", " ")
        while ($newValue.Contains("  ")) {
            $newValue = $newValue.Replace("  ", " ")
        }
        $newValue = $newValue.ToUpper()
    }

    if (-not $newValue.Equals($bValue)) {
        $bCell.Value = $newValue
    }
}
